$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Clear the old "champion" rows (1 & 2), keeping their existing styles ---
$ws.Rows.Item(1).ClearContents()
$ws.Rows.Item(2).ClearContents()
$ws.Rows.Item(1).RowHeight = 12.8
$ws.Rows.Item(2).RowHeight = 12.8

# --- 2. Drop the two old mailto hyperlinks attached to D1/D2 ---
# (deleting while enumerating only removes one hyperlink per pass on this
#  host, so run the sweep twice to guarantee the collection is empty)
foreach ($h in $ws.Hyperlinks) { $h.Delete() }
foreach ($h in $ws.Hyperlinks) { $h.Delete() }

# --- 3. Write the new "Best Singer" notification row (row 4) ---
$ws.Range("A4").Value = "10：00"
$ws.Range("A4").Characters(1, 2).Font.Name = "Arial"
$ws.Range("A4").Characters(3, 1).Font.Name = "Noto Sans CJK SC Regular"
$ws.Range("A4").Characters(4, 2).Font.Name = "Arial"

$ws.Range("A4").Style = "Normal"

$ws.Range("B4").Value = "Zkshadow"
$ws.Range("B4").Style = "Normal"

$ws.Range("C4").Value = "liu676785882@gmail.com"
$ws.Range("C4").Font.Name = "Arial"
$ws.Range("C4").Font.Color = 16711680

$ws.Range("D4").Value = "zkshadow"
$ws.Range("D4").Style = "Normal"

$ws.Range("E4").Value = "2018-11-3 11:33:19"
$ws.Range("E4").Style = "Normal"

# Re-create the hyperlink, now anchored on the email cell only (C4)
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:liu676785882@gmail.com", "", "", "liu676785882@gmail.com")

# --- 4. Row height + selection to match the finished layout ---
$ws.Rows.Item(4).RowHeight = 17.2
$ws.Range("A4").Select()
